# Auto-generated edit script applying numeric updates to Sheets/Famfrit_Profits data
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1999.5
$ws.Range("I53").Value = 1749
$ws.Range("J53").Value = 2250
$ws.Range("K53").Value = 1749
$ws.Range("L53").Value = 2250
$ws.Range("M53").Value = -1112
$ws.Range("N53").Value = -3524
$ws.Range("H86").Value = 2441.0967
$ws.Range("I86").Value = 1523
$ws.Range("J86").Value = 3301.8125
$ws.Range("K86").Value = 1523
$ws.Range("L86").Value = 3301.8125
$ws.Range("M86").Value = -400
$ws.Range("N86").Value = -5547.8125
$ws.Range("H89").Value = 2441.0967
$ws.Range("I89").Value = 1523
$ws.Range("J89").Value = 3301.8125
$ws.Range("K89").Value = 7615
$ws.Range("L89").Value = 16509.0625
$ws.Range("M89").Value = -1999
$ws.Range("N89").Value = -27741.0625
$ws.Range("H106").Value = 1022.0714
$ws.Range("I106").Value = 1043
$ws.Range("K106").Value = 1043
$ws.Range("M106").Value = -412
$ws.Range("H111").Value = 2681.65
$ws.Range("I111").Value = 2646.3333
$ws.Range("J111").Value = 2999.5
$ws.Range("K111").Value = 7938.999899999999
$ws.Range("L111").Value = 8998.5
$ws.Range("M111").Value = -4871.999899999999
$ws.Range("N111").Value = -15132.5
$ws.Range("H116").Value = 2104.2
$ws.Range("I116").Value = 2116
$ws.Range("K116").Value = 2116
$ws.Range("M116").Value = 1326
$ws.Range("H132").Value = 5057.3335
$ws.Range("I132").Value = 5610.7144
$ws.Range("K132").Value = 16832.1432
$ws.Range("M132").Value = -14302.1432
$ws.Range("H138").Value = 6806727.5
$ws.Range("I138").Value = 983.0454999999999
$ws.Range("J138").Value = 12352149
$ws.Range("K138").Value = 2949.1365
$ws.Range("L138").Value = 37056447
$ws.Range("M138").Value = 2190.8635
$ws.Range("N138").Value = -37066727
$ws.Range("H141").Value = 2357.2856
$ws.Range("I141").Value = 2232.6667
$ws.Range("K141").Value = 6698.000100000001
$ws.Range("M141").Value = -1518.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 7082.6665
$ws.Range("I36").Value = 7082.6665
$ws.Range("K36").Value = 7082.6665
$ws.Range("M36").Value = -6736.6665
$ws.Range("H74").Value = 28605004
$ws.Range("I74").Value = 34522730
$ws.Range("J74").Value = 2667.6667
$ws.Range("K74").Value = 34522730
$ws.Range("L74").Value = 2667.6667
$ws.Range("M74").Value = -34521856
$ws.Range("N74").Value = -4415.6667
$ws.Range("H77").Value = 28605004
$ws.Range("I77").Value = 34522730
$ws.Range("J77").Value = 2667.6667
$ws.Range("K77").Value = 172613650
$ws.Range("L77").Value = 13338.3335
$ws.Range("M77").Value = -172609282
$ws.Range("N77").Value = -22074.3335
$ws.Range("H102").Value = 170086.08
$ws.Range("I102").Value = 252732.38
$ws.Range("J102").Value = 4793.5
$ws.Range("K102").Value = 252732.38
$ws.Range("L102").Value = 4793.5
$ws.Range("M102").Value = -251110.38
$ws.Range("N102").Value = -8037.5
$ws.Range("H122").Value = 3710.3684
$ws.Range("I122").Value = 3173.818
$ws.Range("J122").Value = 4448.125
$ws.Range("K122").Value = 9521.454000000002
$ws.Range("L122").Value = 13344.375
$ws.Range("M122").Value = -7071.454000000002
$ws.Range("N122").Value = -18244.375
$ws.Range("H132").Value = 30354094
$ws.Range("I132").Value = 9599.655000000001
$ws.Range("K132").Value = 28798.965
$ws.Range("M132").Value = -26268.965

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4425.3477
$ws.Range("I99").Value = 3686.9333
$ws.Range("J99").Value = 5809.875
$ws.Range("K99").Value = 3686.9333
$ws.Range("L99").Value = 5809.875
$ws.Range("M99").Value = -2188.9333
$ws.Range("N99").Value = -8805.875
$ws.Range("H134").Value = 2486.1482
$ws.Range("I134").Value = 2255
$ws.Range("J134").Value = 2775.0833
$ws.Range("K134").Value = 6765
$ws.Range("L134").Value = 8325.249899999999
$ws.Range("M134").Value = -4230
$ws.Range("N134").Value = -13395.2499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 30000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 30000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 30000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -30226
$ws.Range("H56").Value = 100
$ws.Range("I56").Value = 100
$ws.Range("K56").Value = 100
$ws.Range("M56").Value = 745
$ws.Range("H58").Value = 3715
$ws.Range("I58").Value = 3707.7856
$ws.Range("K58").Value = 3707.7856
$ws.Range("M58").Value = -3504.7856
$ws.Range("H99").Value = 22544
$ws.Range("I99").Value = 23654.285
$ws.Range("J99").Value = 7000
$ws.Range("K99").Value = 23654.285
$ws.Range("L99").Value = 7000
$ws.Range("M99").Value = -22156.285
$ws.Range("N99").Value = -9996
$ws.Range("H126").Value = 22544
$ws.Range("I126").Value = 23654.285
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 70962.855
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -68492.855
$ws.Range("N126").Value = -25940
$ws.Range("H134").Value = 2244.5
$ws.Range("I134").Value = 2021.1111
$ws.Range("K134").Value = 6063.3333
$ws.Range("M134").Value = -3528.3333
$ws.Range("H136").Value = 3715
$ws.Range("I136").Value = 3707.7856
$ws.Range("K136").Value = 11123.3568
$ws.Range("M136").Value = -8573.356800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 876.5
$ws.Range("I5").Value = 927
$ws.Range("J5").Value = 725
$ws.Range("K5").Value = 2781
$ws.Range("L5").Value = 2175
$ws.Range("M5").Value = -2669
$ws.Range("N5").Value = -2399
$ws.Range("H56").Value = 34030.168
$ws.Range("I56").Value = 34030.168
$ws.Range("K56").Value = 34030.168
$ws.Range("M56").Value = -33500.168
$ws.Range("H88").Value = 9316
$ws.Range("J88").Value = 9316
$ws.Range("L88").Value = 27948
$ws.Range("N88").Value = -28804
$ws.Range("H91").Value = 9316
$ws.Range("J91").Value = 9316
$ws.Range("L91").Value = 27948
$ws.Range("N91").Value = -30912
$ws.Range("H131").Value = 1668.2935
$ws.Range("J131").Value = 1720.8795
$ws.Range("L131").Value = 5162.6385
$ws.Range("N131").Value = -15242.6385
$ws.Range("H133").Value = 4200.6
$ws.Range("I133").Value = 4200.6
$ws.Range("K133").Value = 12601.8
$ws.Range("M133").Value = -7541.800000000001
$ws.Range("H134").Value = 3777.3333
$ws.Range("I134").Value = 1226.7273
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 3680.1819
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = 1389.8181
$ws.Range("N134").Value = -55140
$ws.Range("H135").Value = 876.5
$ws.Range("I135").Value = 927
$ws.Range("J135").Value = 725
$ws.Range("K135").Value = 8343
$ws.Range("L135").Value = 6525
$ws.Range("M135").Value = -5808
$ws.Range("N135").Value = -11595

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 175000
$ws.Range("J100").Value = 175000
$ws.Range("L100").Value = 175000
$ws.Range("N100").Value = -177164
$ws.Range("H113").Value = 3896.5
$ws.Range("I113").Value = 2666.6667
$ws.Range("J113").Value = 4423.5713
$ws.Range("K113").Value = 2666.6667
$ws.Range("L113").Value = 4423.5713
$ws.Range("M113").Value = -496.6667000000002
$ws.Range("N113").Value = -8763.5713
$ws.Range("H122").Value = 3110.6924
$ws.Range("I122").Value = 2893.2222
$ws.Range("K122").Value = 8679.6666
$ws.Range("M122").Value = -6229.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1349.4286
$ws.Range("J22").Value = 1449.3334
$ws.Range("L22").Value = 1449.3334
$ws.Range("N22").Value = -2039.3334
$ws.Range("H27").Value = 1349.4286
$ws.Range("J27").Value = 1449.3334
$ws.Range("L27").Value = 1449.3334
$ws.Range("N27").Value = -1663.3334
$ws.Range("H46").Value = 1264.9688
$ws.Range("I46").Value = 494.34784
$ws.Range("J46").Value = 3234.3333
$ws.Range("K46").Value = 494.34784
$ws.Range("L46").Value = 3234.3333
$ws.Range("M46").Value = -306.34784
$ws.Range("N46").Value = -3610.3333
$ws.Range("H61").Value = 3223.4119
$ws.Range("I61").Value = 2300.2727
$ws.Range("J61").Value = 4915.8335
$ws.Range("K61").Value = 2300.2727
$ws.Range("L61").Value = 4915.8335
$ws.Range("M61").Value = -2098.2727
$ws.Range("N61").Value = -5319.8335
$ws.Range("H113").Value = 3223.4119
$ws.Range("I113").Value = 2300.2727
$ws.Range("J113").Value = 4915.8335
$ws.Range("K113").Value = 2300.2727
$ws.Range("L113").Value = 4915.8335
$ws.Range("M113").Value = -130.2727
$ws.Range("N113").Value = -9255.833500000001
$ws.Range("H136").Value = 2278.9697
$ws.Range("I136").Value = 1068.1305
$ws.Range("K136").Value = 3204.3915
$ws.Range("M136").Value = -654.3914999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1293.3636
$ws.Range("I81").Value = 1279
$ws.Range("J81").Value = 1318.5
$ws.Range("K81").Value = 2558
$ws.Range("L81").Value = 2637
$ws.Range("M81").Value = -1497
$ws.Range("N81").Value = -4759
$ws.Range("H84").Value = 1293.3636
$ws.Range("I84").Value = 1279
$ws.Range("J84").Value = 1318.5
$ws.Range("K84").Value = 12790
$ws.Range("L84").Value = 13185
$ws.Range("M84").Value = -7486
$ws.Range("N84").Value = -23793
$ws.Range("H112").Value = 10955.5
$ws.Range("J112").Value = 10955.5
$ws.Range("L112").Value = 10955.5
$ws.Range("N112").Value = -13909.5
$ws.Range("H132").Value = 3684
$ws.Range("I132").Value = 3841.8823
$ws.Range("K132").Value = 11525.6469
$ws.Range("M132").Value = -8995.6469
$ws.Range("H136").Value = 4677.4707
$ws.Range("I136").Value = 1200.1428
$ws.Range("J136").Value = 7111.6
$ws.Range("K136").Value = 3600.4284
$ws.Range("L136").Value = 21334.8
$ws.Range("M136").Value = -1050.4284
$ws.Range("N136").Value = -26434.8

